$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells stay as Text, matching the source sheet,
# instead of being auto-coerced to numbers by Excel on assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "36.406.80"
$ws.Cells.Item(2, 5).Value = "  -0.13%  "
$ws.Cells.Item(3, 4).Value = "1.949.02"
$ws.Cells.Item(3, 5).Value = "  -1.81%  "
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$ws.Cells.Item(5, 4).Value = "243.26"
$ws.Cells.Item(5, 5).Value = "  -0.68%  "
$ws.Cells.Item(6, 4).Value = "0.615"
$ws.Cells.Item(6, 5).Value = "  -0.78%  "
$ws.Cells.Item(7, 5).Value = "  -0.01%  "
$ws.Cells.Item(8, 4).Value = "57.22"
$ws.Cells.Item(8, 5).Value = "  -3.00%  "
$ws.Cells.Item(9, 4).Value = "0.364"
$ws.Cells.Item(9, 5).Value = "  -3.57%  "
$ws.Cells.Item(10, 4).Value = "0.0852"
$ws.Cells.Item(10, 5).Value = "  +3.37%  "
$ws.Cells.Item(11, 5).Value = "  +0.13%  "
$ws.Cells.Item(12, 4).Value = "2.237.99"
$ws.Cells.Item(12, 5).Value = "  -1.66%  "
$ws.Cells.Item(13, 4).Value = "0.819"
$ws.Cells.Item(13, 5).Value = "  -5.53%  "
$ws.Cells.Item(14, 4).Value = "21.25"
$ws.Cells.Item(14, 5).Value = "  -12.28%  "
$ws.Cells.Item(15, 4).Value = "13.55"
$ws.Cells.Item(15, 5).Value = "  -3.83%  "
$ws.Cells.Item(16, 4).Value = "5.20"
$ws.Cells.Item(16, 5).Value = "  -5.12%  "
$ws.Cells.Item(17, 4).Value = "1.960.06"
$ws.Cells.Item(17, 5).Value = "  -1.07%  "
$ws.Cells.Item(18, 4).Value = "36.347.28"
$ws.Cells.Item(18, 5).Value = "  -0.14%  "
$ws.Cells.Item(19, 4).Value = "0.0₃0879"
$ws.Cells.Item(19, 5).Value = "  +1.40%  "
$ws.Cells.Item(20, 4).Value = "69.64"
$ws.Cells.Item(20, 5).Value = "  -1.72%  "
$ws.Cells.Item(21, 4).Value = "229.69"
$ws.Cells.Item(21, 5).Value = "  -2.24%  "
$ws.Cells.Item(22, 4).Value = "5.05"
$ws.Cells.Item(22, 5).Value = "  -5.58%  "
$ws.Cells.Item(23, 5).Value = "  +0.00%  "
$ws.Cells.Item(24, 4).Value = "2.41"
$ws.Cells.Item(24, 5).Value = "  -8.37%  "
$ws.Cells.Item(25, 4).Value = "2.29"
$ws.Cells.Item(25, 5).Value = "  -0.73%  "
$ws.Cells.Item(26, 4).Value = "9.24"
$ws.Cells.Item(26, 5).Value = "  -9.27%  "
$ws.Cells.Item(27, 2).Value = "Kaspa"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(27, 4).Value = "0.138"
$ws.Cells.Item(27, 5).Value = "  +10.19%  "
$ws.Cells.Item(28, 2).Value = "Monero"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(28, 4).Value = "161.20"
$ws.Cells.Item(28, 5).Value = "  -0.49%  "
$ws.Cells.Item(29, 4).Value = "19.34"
$ws.Cells.Item(29, 5).Value = "  -2.93%  "
$ws.Cells.Item(30, 5).Value = "  -1.95%  "
$ws.Cells.Item(31, 4).Value = "1.15"
$ws.Cells.Item(31, 5).Value = "  -4.31%  "
$ws.Cells.Item(32, 4).Value = "4.63"
$ws.Cells.Item(32, 5).Value = "  -6.05%  "
$ws.Cells.Item(33, 4).Value = "0.0648"
$ws.Cells.Item(33, 5).Value = "  +2.05%  "
$ws.Cells.Item(34, 4).Value = "4.28"
$ws.Cells.Item(34, 5).Value = "  -3.87%  "
$ws.Cells.Item(35, 4).Value = "6.18"
$ws.Cells.Item(35, 5).Value = "  -2.48%  "
$ws.Cells.Item(36, 5).Value = "  +0.04%  "
$ws.Cells.Item(37, 5).Value = "  +1.34%  "
$ws.Cells.Item(38, 4).Value = "2.15"
$ws.Cells.Item(38, 5).Value = "  -5.93%  "
$ws.Cells.Item(39, 5).Value = "  -1.64%  "
$ws.Cells.Item(40, 4).Value = "0.0978"
$ws.Cells.Item(40, 5).Value = "  +0.89%  "
$ws.Cells.Item(41, 5).Value = "  +0.32%  "
$ws.Cells.Item(42, 2).Value = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(42, 4).Value = "1.17"
$ws.Cells.Item(42, 5).Value = "  -7.13%  "
$ws.Cells.Item(43, 2).Value = "VeChain"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(43, 4).Value = "0.0212"
$ws.Cells.Item(43, 5).Value = "  -1.53%  "
$ws.Cells.Item(44, 2).Value = "Maker"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(44, 4).Value = "1.356.06"
$ws.Cells.Item(44, 5).Value = "  -1.56%  "
$ws.Cells.Item(45, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(45, 4).Value = "15.64"
$ws.Cells.Item(45, 5).Value = "  -4.57%  "
$ws.Cells.Item(46, 5).Value = "  -6.66%  "
$ws.Cells.Item(47, 4).Value = "87.47"
$ws.Cells.Item(47, 5).Value = "  -5.92%  "
$ws.Cells.Item(48, 4).Value = "7.13"
$ws.Cells.Item(48, 5).Value = "  -6.65%  "
$ws.Cells.Item(49, 4).Value = "2.83"
$ws.Cells.Item(49, 5).Value = "  -1.12%  "
$ws.Cells.Item(50, 4).Value = "44.76"
$ws.Cells.Item(50, 5).Value = "  -1.54%  "
$ws.Cells.Item(51, 4).Value = "2.128.34"
$ws.Cells.Item(51, 5).Value = "  -1.93%  "
